$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Move the existing "Deadlines" values (currently in column F) over to the
#    new column I, freeing up F for the new "Anton" presence column, and
#    insert two new columns G ("Samuel") and H ("Frida") in between.
# ---------------------------------------------------------------------------
$deadlineRows = @(1, 3, 10, 25, 30, 38, 39, 45, 52, 53, 66)
foreach ($r in $deadlineRows) {
    $val = $ws.Range("F$r").Value()
    $ws.Range("I$r").Value = $val
    $ws.Range("F$r").Value = $null
}

# ---------------------------------------------------------------------------
# 2. Header row: add the three new "present" columns between the Gatherings
#    column (E) and the Deadlines column (now I).
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "Anton"
$ws.Range("G1").Value = "Samuel"
$ws.Range("H1").Value = "Frida"

# ---------------------------------------------------------------------------
# 3. Fill in the possible dates Anton (F), Samuel (G) and Frida (H) can do
#    the laboration - "X" means available, "Kanske" means maybe.
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = "X"
$ws.Range("G3").Value = "X"
$ws.Range("H3").Value = "X"

$ws.Range("F9").Value = "X"

$ws.Range("F10").Value = "X"

$ws.Range("G11").Value = "X"

$ws.Range("F16").Value = "X"
$ws.Range("G16").Value = "X"

$ws.Range("F17").Value = "X"

$ws.Range("F25").Value = "X"

$ws.Range("F30").Value = "X"

$ws.Range("F38").Value = "X"

$ws.Range("F39").Value = "Kanske"

$ws.Range("F44").Value = "Kanske"

$ws.Range("F45").Value = "Kanske"

$ws.Range("F51").Value = "X"

$ws.Range("F52").Value = "X"

$ws.Range("F53").Value = "Kanske"

$ws.Range("F59").Value = "X"

# ---------------------------------------------------------------------------
# 4. Column widths - shrink the former "Deadlines" column (F) now that it
#    only needs to fit "Anton"/"X" and size the two new columns, while I
#    (the relocated Deadlines column) keeps roughly the old F width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 4.833333333333333
$ws.Columns.Item(3).ColumnWidth = 5.5
$ws.Columns.Item(4).ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666
$ws.Columns.Item(6).ColumnWidth = 5.5
$ws.Columns.Item(7).ColumnWidth = 5.666666666666667
$ws.Columns.Item(8).ColumnWidth = 4
$ws.Columns.Item(9).ColumnWidth = 11.5

# ---------------------------------------------------------------------------
# 5. Freeze the header row and set the active selection/cell.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F5").Select() | Out-Null
